# Adds a "Class" label (flagged as a misspelling by the es-CO proofing
# language, same as the author's manual edit) in front of each "Void"
# keyword used as a (made-up) return-type label on slides 13-15, and
# tidies up the run that follows so the word-break still reads correctly.
#
# Two flavours show up in the source deck:
#   1) "Void " (trailing space already inside the Void run) -> splits
#      naturally into "Class" + " " when we overwrite just the word.
#   2) "Void" (no trailing space, the next run starts with a leading
#      space) -> usually just becomes "Class" in place; for the single
#      "jugador();" paragraph the leading space instead moves into a
#      dedicated middle run (matching the source edit) and the next
#      run's own leading space is trimmed.

$p = $ppt.ActivePresentation

function Replace-VoidSimple($textRange, $afterPos) {
    # "Void " (with trailing space) -> "Class" (new run) + " " (old run).
    $found = $textRange.Find("Void", $afterPos, $true)
    $found.Text = "Class"
    return $found.Start + 5
}

function Replace-VoidInPlace($textRange, $afterPos) {
    # "Void" (no trailing space) -> "Class", same single run.
    $found = $textRange.Find("Void", $afterPos, $true)
    $found.Text = "Class"
    return $found.Start + 5
}

function Replace-VoidSplitAcrossRun($textRange, $afterPos) {
    # "Void" immediately followed by a run that begins with a space
    # (" jugador();") -> "Class" + " " as two distinct runs, and the
    # following run loses its now-redundant leading space.
    $found = $textRange.Find("Void", $afterPos, $true)
    $start = $found.Start

    # Pull "Void" plus the leading space of the following run (5 chars)
    # into one range and retype it as "Class " - this also eats the
    # next run's leading space.
    $crossRun = $textRange.Characters($start, 5)
    $crossRun.Text = "Class "

    # Re-touch just the "Class" portion so it splits off into its own
    # run, leaving the trailing space as a separate run (matching the
    # source file's run layout).
    $classOnly = $textRange.Characters($start, 5)
    $classOnly.Text = "Class"

    return $start + 6
}

# ---- Slide 13 : "CuadroTexto 3" (shape 2) ----
$s13 = $p.Slides.Item(13)
$tr13 = $s13.Shapes.Item(2).TextFrame.TextRange

$pos = 0
$pos = Replace-VoidSimple $tr13 $pos          # Void escenario();
$pos = Replace-VoidSimple $tr13 $pos          # Void puente();
$pos = Replace-VoidSplitAcrossRun $tr13 $pos  # Void jugador();
$pos = Replace-VoidInPlace $tr13 $pos         # Void enemigos();
$pos = Replace-VoidInPlace $tr13 $pos         # Void objeto();

# ---- Slide 14 : "CuadroTexto 3" (shape 1) ----
$s14 = $p.Slides.Item(14)
$tr14 = $s14.Shapes.Item(1).TextFrame.TextRange

$pos = 0
$pos = Replace-VoidSimple $tr14 $pos          # Void juego();
$found = $tr14.Find("Void", $pos, $true)      # Void actualizar_juego(); (untouched)
$pos = $found.Start + $found.Length
$pos = Replace-VoidInPlace $tr14 $pos         # Void focus();

# ---- Slide 15 : "CuadroTexto 3" (shape 1) ----
$s15 = $p.Slides.Item(15)
$tr15 = $s15.Shapes.Item(1).TextFrame.TextRange

$pos = 0
$pos = Replace-VoidSimple $tr15 $pos          # Void mini_juego();
